$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.134653
$ws.Range("H2").Value = 0.403959
$ws.Range("I2").Value = 0.0184956077859211
$ws.Range("J2").Value = 0.01849560778592109
$ws.Range("M2").Value = 209.26237
$ws.Range("N2").Value = 627.78711
$ws.Range("O2").Value = 0.8127157202241573
$ws.Range("P2").Value = 0.8127157202241573
$ws.Range("Q2").Value = 28.17780590761
$ws.Range("R2").Value = 253.60025316849
$ws.Range("S2").Value = 0.01503167120271839
$ws.Range("T2").Value = 0.01503167120271839

# Row 3
$ws.Range("G3").Value = 0.134653
$ws.Range("H3").Value = 0.403959
$ws.Range("I3").Value = 0.0184956077859211
$ws.Range("J3").Value = 0.01849560778592109
$ws.Range("M3").Value = 0.9848756666666668
$ws.Range("N3").Value = 2.954627
$ws.Range("O3").Value = 0.003824977881910862
$ws.Range("P3").Value = 0.003824977881910862
$ws.Range("Q3").Value = 0.1326164631436667
$ws.Range("R3").Value = 1.193548168293
$ws.Range("S3").Value = 0.00007074529069364652
$ws.Range("T3").Value = 0.00007074529069364651

# Row 4
$ws.Range("G4").Value = 0.134653
$ws.Range("H4").Value = 0.403959
$ws.Range("I4").Value = 0.0184956077859211
$ws.Range("J4").Value = 0.01849560778592109
$ws.Range("M4").Value = 1.763846666666667
$ws.Range("N4").Value = 5.291539999999999
$ws.Range("O4").Value = 0.006850280411451801
$ws.Range("P4").Value = 0.006850280411451801
$ws.Range("Q4").Value = 0.2375072452066666
$ws.Range("R4").Value = 2.13756520686
$ws.Range("S4").Value = 0.0001267000997137907
$ws.Range("T4").Value = 0.0001267000997137907

# Row 5
$ws.Range("G5").Value = 0.134653
$ws.Range("H5").Value = 0.403959
$ws.Range("I5").Value = 0.0184956077859211
$ws.Range("J5").Value = 0.01849560778592109
$ws.Range("M5").Value = 45.474231
$ws.Range("N5").Value = 136.422693
$ws.Range("O5").Value = 0.1766090214824801
$ws.Range("P5").Value = 0.1766090214824801
$ws.Range("Q5").Value = 6.123241626842999
$ws.Range("R5").Value = 55.10917464158699
$ws.Range("S5").Value = 0.003266491192795265
$ws.Range("T5").Value = 0.003266491192795264

# Row 6
$ws.Range("I6").Value = 0.5431242536047317
$ws.Range("J6").Value = 0.5431242536047317
$ws.Range("M6").Value = 209.26237
$ws.Range("N6").Value = 627.78711
$ws.Range("O6").Value = 0.8127157202241573
$ws.Range("P6").Value = 0.8127157202241573
$ws.Range("Q6").Value = 827.4423841015467
$ws.Range("R6").Value = 7446.981456913921
$ws.Range("S6").Value = 0.4414056189395774
$ws.Range("T6").Value = 0.4414056189395774

# Row 7
$ws.Range("I7").Value = 0.5431242536047317
$ws.Range("J7").Value = 0.5431242536047317
$ws.Range("M7").Value = 0.9848756666666668
$ws.Range("N7").Value = 2.954627
$ws.Range("O7").Value = 0.003824977881910862
$ws.Range("P7").Value = 0.003824977881910862
$ws.Range("Q7").Value = 3.894287681393779
$ws.Range("R7").Value = 35.04858913254401
$ws.Range("S7").Value = 0.002077438257167445
$ws.Range("T7").Value = 0.002077438257167445

# Row 8
$ws.Range("I8").Value = 0.5431242536047317
$ws.Range("J8").Value = 0.5431242536047317
$ws.Range("M8").Value = 1.763846666666667
$ws.Range("N8").Value = 5.291539999999999
$ws.Range("O8").Value = 0.006850280411451801
$ws.Range("P8").Value = 0.006850280411451801
$ws.Range("Q8").Value = 6.974409642097778
$ws.Range("R8").Value = 62.76968677888
$ws.Range("S8").Value = 0.003720553435452874
$ws.Range("T8").Value = 0.003720553435452874

# Row 9
$ws.Range("I9").Value = 0.5431242536047317
$ws.Range("J9").Value = 0.5431242536047317
$ws.Range("M9").Value = 45.474231
$ws.Range("N9").Value = 136.422693
$ws.Range("O9").Value = 0.1766090214824801
$ws.Range("P9").Value = 0.1766090214824801
$ws.Range("Q9").Value = 179.809232370944
$ws.Range("R9").Value = 1618.283091338496
$ws.Range("S9").Value = 0.09592064297253403
$ws.Range("T9").Value = 0.09592064297253403

# Row 10
$ws.Range("G10").Value = 3.191525333333333
$ws.Range("H10").Value = 9.574576
$ws.Range("I10").Value = 0.4383801386093472
$ws.Range("J10").Value = 0.4383801386093472
$ws.Range("M10").Value = 209.26237
$ws.Range("N10").Value = 627.78711
$ws.Range("O10").Value = 0.8127157202241573
$ws.Range("P10").Value = 0.8127157202241573
$ws.Range("Q10").Value = 667.8661551683733
$ws.Range("R10").Value = 6010.79539651536
$ws.Range("S10").Value = 0.3562784300818615
$ws.Range("T10").Value = 0.3562784300818615

# Row 11
$ws.Range("G11").Value = 3.191525333333333
$ws.Range("H11").Value = 9.574576
$ws.Range("I11").Value = 0.4383801386093472
$ws.Range("J11").Value = 0.4383801386093472
$ws.Range("M11").Value = 0.9848756666666668
$ws.Range("N11").Value = 2.954627
$ws.Range("O11").Value = 0.003824977881910862
$ws.Range("P11").Value = 0.003824977881910862
$ws.Range("Q11").Value = 3.143255640350223
$ws.Range("R11").Value = 28.289300763152
$ws.Range("S11").Value = 0.001676794334049771
$ws.Range("T11").Value = 0.001676794334049771

# Row 12
$ws.Range("G12").Value = 3.191525333333333
$ws.Range("H12").Value = 9.574576
$ws.Range("I12").Value = 0.4383801386093472
$ws.Range("J12").Value = 0.4383801386093472
$ws.Range("M12").Value = 1.763846666666667
$ws.Range("N12").Value = 5.291539999999999
$ws.Range("O12").Value = 0.006850280411451801
$ws.Range("P12").Value = 0.006850280411451801
$ws.Range("Q12").Value = 5.629361320782222
$ws.Range("R12").Value = 50.66425188704
$ws.Range("S12").Value = 0.003003026876285136
$ws.Range("T12").Value = 0.003003026876285136

# Row 13
$ws.Range("G13").Value = 3.191525333333333
$ws.Range("H13").Value = 9.574576
$ws.Range("I13").Value = 0.4383801386093472
$ws.Range("J13").Value = 0.4383801386093472
$ws.Range("M13").Value = 45.474231
$ws.Range("N13").Value = 136.422693
$ws.Range("O13").Value = 0.1766090214824801
$ws.Range("P13").Value = 0.1766090214824801
$ws.Range("Q13").Value = 145.132160250352
$ws.Range("R13").Value = 1306.189442253168
$ws.Range("S13").Value = 0.07742188731715079
$ws.Range("T13").Value = 0.07742188731715079

